$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14
$ws.Range("B3").Value = 212
$ws.Range("B4").Value = 78
$ws.Range("B5").Value = 495
$ws.Range("B6").Value = 14
$ws.Range("B8").Value = 16
$ws.Range("B9").Value = 9
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 12
$ws.Range("B15").Value = 669
